# Auto-generated edit script: apply 2024-08-30 crime data update
# Updates the 2024 (column K) counts across the Citywide Totals,
# By Neighborhood, and individual neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5282
$ws.Range("K3").Value = 5448
$ws.Range("K4").Value = 1130
$ws.Range("K5").Value = 388
$ws.Range("K6").Value = 6064
$ws.Range("K7").Value = 18312

$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("K5").Value = 8
$ws.Range("K6").Value = 17

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 338
$ws.Range("K3").Value = 369
$ws.Range("K7").Value = 1232

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 136
$ws.Range("K7").Value = 406

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 215
$ws.Range("K3").Value = 290
$ws.Range("K6").Value = 228
$ws.Range("K7").Value = 785

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 177
$ws.Range("K6").Value = 180
$ws.Range("K7").Value = 617

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 133
$ws.Range("K6").Value = 154
$ws.Range("K7").Value = 416

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 128
$ws.Range("K7").Value = 309

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 161
$ws.Range("K7").Value = 541
$ws.Range("K8").Value = 1232
$ws.Range("K10").Value = 101
$ws.Range("K11").Value = 350
$ws.Range("K14").Value = 96
$ws.Range("K15").Value = 187
$ws.Range("K19").Value = 538
$ws.Range("K20").Value = 423
$ws.Range("K23").Value = 193
$ws.Range("K27").Value = 175
$ws.Range("K29").Value = 980
$ws.Range("K30").Value = 71
$ws.Range("K33").Value = 785
$ws.Range("K34").Value = 104
$ws.Range("K37").Value = 617
$ws.Range("K38").Value = 17
$ws.Range("K42").Value = 677
$ws.Range("K43").Value = 163
$ws.Range("K47").Value = 124
$ws.Range("K48").Value = 230
$ws.Range("K52").Value = 477
$ws.Range("K54").Value = 360
$ws.Range("K60").Value = 114
$ws.Range("K63").Value = 52
$ws.Range("K64").Value = 117
$ws.Range("K65").Value = 416
$ws.Range("K66").Value = 60
$ws.Range("K67").Value = 694
$ws.Range("K73").Value = 157
$ws.Range("K76").Value = 253
$ws.Range("K78").Value = 210
$ws.Range("K79").Value = 450
$ws.Range("K83").Value = 406
$ws.Range("K85").Value = 864
$ws.Range("K89").Value = 265
$ws.Range("K90").Value = 165
$ws.Range("K91").Value = 201
$ws.Range("K94").Value = 241
$ws.Range("K97").Value = 146
$ws.Range("K98").Value = 86
$ws.Range("K99").Value = 309
$ws.Range("K101").Value = 18312

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 199
$ws.Range("K3").Value = 247
$ws.Range("K7").Value = 694

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 90
$ws.Range("K6").Value = 192
$ws.Range("K7").Value = 360

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 280
$ws.Range("K3").Value = 354
$ws.Range("K6").Value = 271
$ws.Range("K7").Value = 980

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 114
$ws.Range("K7").Value = 230

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 158
$ws.Range("K3").Value = 169
$ws.Range("K7").Value = 538

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 55
$ws.Range("K7").Value = 253

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 179
$ws.Range("K6").Value = 256
$ws.Range("K7").Value = 677

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K3").Value = 17
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 47
$ws.Range("K7").Value = 210

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 53
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 193

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 50
$ws.Range("K3").Value = 96
$ws.Range("K7").Value = 201

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K6").Value = 110
$ws.Range("K7").Value = 450

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 117

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 139
$ws.Range("K5").Value = 8
$ws.Range("K7").Value = 423

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 184
$ws.Range("K6").Value = 144
$ws.Range("K7").Value = 541

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 104
$ws.Range("K7").Value = 241

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 124

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 65
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K3").Value = 13
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K4").Value = 21
$ws.Range("K7").Value = 350

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 51
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K4").Value = 31
$ws.Range("K7").Value = 265

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 48
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 175

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 59
$ws.Range("K7").Value = 165

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K4").Value = 23
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 163

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 289
$ws.Range("K3").Value = 292
$ws.Range("K7").Value = 864

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 127
$ws.Range("K3").Value = 134
$ws.Range("K7").Value = 477
